$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FL1013"
$ws.Range("B2").Value = "app"
$ws.Range("C2").Value = "business"
$ws.Range("D2").Value = 8.0
$ws.Range("E2").Value = 5501.03
$ws.Range("F2").Value = 366.46
$ws.Range("G2").Value = 850.0
$ws.Range("H2").Value = 6717.49
$ws.Range("I2").Value = 839.69
$ws.Range("J2").Value = 335.87

$ws.Range("A3").Value = "FL1013"
$ws.Range("B3").Value = "app"
$ws.Range("C3").Value = "economy"
$ws.Range("D3").Value = 106.0
$ws.Range("E3").Value = 47213.44
$ws.Range("F3").Value = 5241.3
$ws.Range("G3").Value = 7300.0
$ws.Range("H3").Value = 59754.74
$ws.Range("I3").Value = 563.72
$ws.Range("J3").Value = 2987.74

$ws.Range("A4").Value = "FL1013"
$ws.Range("B4").Value = "app"
$ws.Range("C4").Value = "premium_economy"
$ws.Range("D4").Value = 6.0
$ws.Range("E4").Value = 3116.38
$ws.Range("F4").Value = 362.24
$ws.Range("G4").Value = 450.0
$ws.Range("H4").Value = 3928.62
$ws.Range("I4").Value = 654.77
$ws.Range("J4").Value = 196.43

$ws.Range("A5").Value = "FL1013"
$ws.Range("B5").Value = "travel_agent"
$ws.Range("C5").Value = "business"
$ws.Range("D5").Value = 10.0
$ws.Range("E5").Value = 6820.37
$ws.Range("F5").Value = 418.1
$ws.Range("G5").Value = 500.0
$ws.Range("H5").Value = 7738.47
$ws.Range("I5").Value = 773.85
$ws.Range("J5").Value = 386.92

$ws.Range("A6").Value = "FL1013"
$ws.Range("B6").Value = "travel_agent"
$ws.Range("C6").Value = "economy"
$ws.Range("D6").Value = 92.0
$ws.Range("E6").Value = 40852.81
$ws.Range("F6").Value = 4767.7
$ws.Range("G6").Value = 6650.0
$ws.Range("H6").Value = 52270.51
$ws.Range("I6").Value = 568.16
$ws.Range("J6").Value = 2613.53

$ws.Range("A7").Value = "FL1013"
$ws.Range("B7").Value = "travel_agent"
$ws.Range("C7").Value = "premium_economy"
$ws.Range("D7").Value = 25.0
$ws.Range("E7").Value = 14391.32
$ws.Range("F7").Value = 1241.59
$ws.Range("G7").Value = 2050.0
$ws.Range("H7").Value = 17682.91
$ws.Range("I7").Value = 707.32
$ws.Range("J7").Value = 884.15

$ws.Range("A8").Value = "FL1013"
$ws.Range("B8").Value = "website"
$ws.Range("C8").Value = "business"
$ws.Range("D8").Value = 7.0
$ws.Range("E8").Value = 4652.36
$ws.Range("F8").Value = 290.12
$ws.Range("G8").Value = 800.0
$ws.Range("H8").Value = 5742.48
$ws.Range("I8").Value = 820.35
$ws.Range("J8").Value = 287.12

$ws.Range("A9").Value = "FL1013"
$ws.Range("B9").Value = "website"
$ws.Range("C9").Value = "economy"
$ws.Range("D9").Value = 85.0
$ws.Range("E9").Value = 37283.53
$ws.Range("F9").Value = 3738.5
$ws.Range("G9").Value = 6750.0
$ws.Range("H9").Value = 47772.03
$ws.Range("I9").Value = 562.02
$ws.Range("J9").Value = 2388.6

$ws.Range("A10").Value = "FL1013"
$ws.Range("B10").Value = "website"
$ws.Range("C10").Value = "premium_economy"
$ws.Range("D10").Value = 25.0
$ws.Range("E10").Value = 13937.32
$ws.Range("F10").Value = 1171.61
$ws.Range("G10").Value = 2150.0
$ws.Range("H10").Value = 17258.93
$ws.Range("I10").Value = 690.36
$ws.Range("J10").Value = 862.95

$ws.Range("A11").Value = "FL1020"
$ws.Range("B11").Value = "app"
$ws.Range("C11").Value = "business"
$ws.Range("D11").Value = 10.0
$ws.Range("E11").Value = 6444.27
$ws.Range("F11").Value = 499.96
$ws.Range("G11").Value = 650.0
$ws.Range("H11").Value = 7594.23
$ws.Range("I11").Value = 759.42
$ws.Range("J11").Value = 379.71

$ws.Range("A12").Value = "FL1020"
$ws.Range("B12").Value = "app"
$ws.Range("C12").Value = "economy"
$ws.Range("D12").Value = 68.0
$ws.Range("E12").Value = 29862.23
$ws.Range("F12").Value = 3511.68
$ws.Range("G12").Value = 4400.0
$ws.Range("H12").Value = 37773.91
$ws.Range("I12").Value = 555.5
$ws.Range("J12").Value = 1888.7

$ws.Range("A13").Value = "FL1020"
$ws.Range("B13").Value = "app"
$ws.Range("C13").Value = "premium_economy"
$ws.Range("D13").Value = 8.0
$ws.Range("E13").Value = 4683.6
$ws.Range("F13").Value = 213.85
$ws.Range("G13").Value = 700.0
$ws.Range("H13").Value = 5597.45
$ws.Range("I13").Value = 699.68
$ws.Range("J13").Value = 279.87

$ws.Range("A14").Value = "FL1020"
$ws.Range("B14").Value = "travel_agent"
$ws.Range("C14").Value = "business"
$ws.Range("D14").Value = 7.0
$ws.Range("E14").Value = 4512.76
$ws.Range("F14").Value = 442.38
$ws.Range("G14").Value = 300.0
$ws.Range("H14").Value = 5255.14
$ws.Range("I14").Value = 750.73
$ws.Range("J14").Value = 262.76

$ws.Range("A15").Value = "FL1020"
$ws.Range("B15").Value = "travel_agent"
$ws.Range("C15").Value = "economy"
$ws.Range("D15").Value = 55.0
$ws.Range("E15").Value = 24292.52
$ws.Range("F15").Value = 2585.68
$ws.Range("G15").Value = 3850.0
$ws.Range("H15").Value = 30728.2
$ws.Range("I15").Value = 558.69
$ws.Range("J15").Value = 1536.41

$ws.Range("A16").Value = "FL1020"
$ws.Range("B16").Value = "travel_agent"
$ws.Range("C16").Value = "premium_economy"
$ws.Range("D16").Value = 13.0
$ws.Range("E16").Value = 7149.17
$ws.Range("F16").Value = 682.17
$ws.Range("G16").Value = 850.0
$ws.Range("H16").Value = 8681.34
$ws.Range("I16").Value = 667.8
$ws.Range("J16").Value = 434.07

$ws.Range("A17").Value = "FL1020"
$ws.Range("B17").Value = "website"
$ws.Range("C17").Value = "business"
$ws.Range("D17").Value = 8.0
$ws.Range("E17").Value = 5468.44
$ws.Range("F17").Value = 384.19
$ws.Range("G17").Value = 950.0
$ws.Range("H17").Value = 6802.63
$ws.Range("I17").Value = 850.33
$ws.Range("J17").Value = 340.13

$ws.Range("A18").Value = "FL1020"
$ws.Range("B18").Value = "website"
$ws.Range("C18").Value = "economy"
$ws.Range("D18").Value = 56.0
$ws.Range("E18").Value = 24377.77
$ws.Range("F18").Value = 3056.29
$ws.Range("G18").Value = 4950.0
$ws.Range("H18").Value = 32384.06
$ws.Range("I18").Value = 578.29
$ws.Range("J18").Value = 1619.2

$ws.Range("A19").Value = "FL1020"
$ws.Range("B19").Value = "website"
$ws.Range("C19").Value = "premium_economy"
$ws.Range("D19").Value = 7.0
$ws.Range("E19").Value = 3739.54
$ws.Range("F19").Value = 452.08
$ws.Range("G19").Value = 650.0
$ws.Range("H19").Value = 4841.62
$ws.Range("I19").Value = 691.66
$ws.Range("J19").Value = 242.08

$ws.Range("A20").Value = "FL1023"
$ws.Range("B20").Value = "app"
$ws.Range("C20").Value = "business"
$ws.Range("D20").Value = 3.0
$ws.Range("E20").Value = 762.79
$ws.Range("F20").Value = 160.32
$ws.Range("G20").Value = 50.0
$ws.Range("H20").Value = 973.11
$ws.Range("I20").Value = 324.37
$ws.Range("J20").Value = 48.66

$ws.Range("A21").Value = "FL1023"
$ws.Range("B21").Value = "app"
$ws.Range("C21").Value = "economy"
$ws.Range("D21").Value = 35.0
$ws.Range("E21").Value = 5774.24
$ws.Range("F21").Value = 1701.91
$ws.Range("G21").Value = 2650.0
$ws.Range("H21").Value = 10126.15
$ws.Range("I21").Value = 289.32
$ws.Range("J21").Value = 506.31

$ws.Range("A22").Value = "FL1023"
$ws.Range("B22").Value = "app"
$ws.Range("C22").Value = "premium_economy"
$ws.Range("D22").Value = 7.0
$ws.Range("E22").Value = 1430.57
$ws.Range("F22").Value = 357.0
$ws.Range("G22").Value = 600.0
$ws.Range("H22").Value = 2387.57
$ws.Range("I22").Value = 341.08
$ws.Range("J22").Value = 119.38

$ws.Range("A23").Value = "FL1023"
$ws.Range("B23").Value = "travel_agent"
$ws.Range("C23").Value = "business"
$ws.Range("D23").Value = 5.0
$ws.Range("E23").Value = 1156.93
$ws.Range("F23").Value = 255.81
$ws.Range("G23").Value = 350.0
$ws.Range("H23").Value = 1762.74
$ws.Range("I23").Value = 352.55
$ws.Range("J23").Value = 88.14

$ws.Range("A24").Value = "FL1023"
$ws.Range("B24").Value = "travel_agent"
$ws.Range("C24").Value = "economy"
$ws.Range("D24").Value = 34.0
$ws.Range("E24").Value = 5410.87
$ws.Range("F24").Value = 1671.94
$ws.Range("G24").Value = 2300.0
$ws.Range("H24").Value = 9382.81
$ws.Range("I24").Value = 275.97
$ws.Range("J24").Value = 469.14

$ws.Range("A25").Value = "FL1023"
$ws.Range("B25").Value = "travel_agent"
$ws.Range("C25").Value = "premium_economy"
$ws.Range("D25").Value = 7.0
$ws.Range("E25").Value = 1366.89
$ws.Range("F25").Value = 303.67
$ws.Range("G25").Value = 800.0
$ws.Range("H25").Value = 2470.56
$ws.Range("I25").Value = 352.94
$ws.Range("J25").Value = 123.53

$ws.Range("A26").Value = "FL1023"
$ws.Range("B26").Value = "website"
$ws.Range("C26").Value = "business"
$ws.Range("D26").Value = 1.0
$ws.Range("E26").Value = 264.61
$ws.Range("F26").Value = 4.94
$ws.Range("G26").Value = 100.0
$ws.Range("H26").Value = 369.55
$ws.Range("I26").Value = 369.55
$ws.Range("J26").Value = 18.48

